$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.911.67"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.814.40"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3656"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07351"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("D12").Value = "1.785.60"
$ws.Range("E12").Value = "  -1.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.362"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07100"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.497"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008694"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("D21").Value = "26.919.52"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.287"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "2.053.47"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.895"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.251"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08903"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7538"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.478"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.902"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05278"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01948"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.967"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.237"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5310"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1652"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.415"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4866"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.660"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06291"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
